# msz - restructuring control processing -> container
# Adds three new rows to the Tabelle1 control table describing the new
# "Abwesenheit anlegen" / "Abwesenheit Sylvester" container controls, widens
# column A to fit the longer control names, and shifts/resizes the
# screenshot picture further down the sheet to make room for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (10-12) -------------------------------------------------
$ws.Range("A10").Value = "Abwesenheit anlegen Page Abwesenheiten"
$ws.Range("B10").Value = "<SET>"
$ws.Range("D10").Value = "Abwesenheit anlegen"
$ws.Range("F10").Value = "<NOP>"

$ws.Range("A11").Value = "Abwesenheit Sylvester checken"
$ws.Range("B11").Value = "<CHK>"
$ws.Range("D11").Value = "Abwesenheit Sylvester checken"
$ws.Range("F11").Value = "<NOP>"

$ws.Range("A12").Value = "Abwesenheit Sylvester löschen"
$ws.Range("B12").Value = "<SET>"
$ws.Range("D12").Value = "Abwesenheit Sylvester löschen"
$ws.Range("F12").Value = "<NOP>"

# --- Column A is now wider to fit the new, longer control names -----------
$ws.Columns.Item(1).ColumnWidth = 46.1666667

# --- Selection moves to the first cell of the freshly added block ---------
$ws.Range("D14").Select()

# --- Move/resize the screenshot picture down to sit below the new rows ----
$shp = $ws.Shapes.Item(1)
$shp.Left = 1.2
$shp.Top = 215.4
$shp.Width = 1169.1038582677165
$shp.Height = 661.4173228346457
